$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 849.5
$ws.Range("J12").Value = 1000.6667
$ws.Range("L12").Value = 1000.6667
$ws.Range("N12").Value = -1340.6667
$ws.Range("H70").Value = 6171
$ws.Range("I70").Value = 14179
$ws.Range("J70").Value = 1900.0667
$ws.Range("K70").Value = 42537
$ws.Range("L70").Value = 5700.2001
$ws.Range("M70").Value = -42267
$ws.Range("N70").Value = -6240.2001
$ws.Range("H73").Value = 6171
$ws.Range("I73").Value = 14179
$ws.Range("J73").Value = 1900.0667
$ws.Range("K73").Value = 42537
$ws.Range("L73").Value = 5700.2001
$ws.Range("M73").Value = -41601
$ws.Range("N73").Value = -7572.2001
$ws.Range("H86").Value = 2288.9092
$ws.Range("J86").Value = 2115
$ws.Range("L86").Value = 2115
$ws.Range("N86").Value = -4361
$ws.Range("H89").Value = 2288.9092
$ws.Range("J89").Value = 2115
$ws.Range("L89").Value = 10575
$ws.Range("N89").Value = -21807
$ws.Range("H106").Value = 11727.272
$ws.Range("I106").Value = 2143.5715
$ws.Range("J106").Value = 28498.75
$ws.Range("K106").Value = 2143.5715
$ws.Range("L106").Value = 28498.75
$ws.Range("M106").Value = -1512.5715
$ws.Range("N106").Value = -29760.75
$ws.Range("H138").Value = 2228.5
$ws.Range("I138").Value = 1445.8
$ws.Range("K138").Value = 4337.4
$ws.Range("M138").Value = 802.6000000000004
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H48").Value = 124998
$ws.Range("J48").Value = 124998
$ws.Range("L48").Value = 124998
$ws.Range("N48").Value = -125766
$ws.Range("H102").Value = 1499.2727
$ws.Range("I102").Value = 1434.1786
$ws.Range("K102").Value = 1434.1786
$ws.Range("M102").Value = 187.8214
$ws.Range("H132").Value = 13043.474
$ws.Range("I132").Value = 14552.4375
$ws.Range("J132").Value = 4995.6665
$ws.Range("K132").Value = 43657.3125
$ws.Range("L132").Value = 14986.9995
$ws.Range("M132").Value = -41127.3125
$ws.Range("N132").Value = -20046.9995
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21049.4
$ws.Range("I82").Value = 11838
$ws.Range("J82").Value = 24997.143
$ws.Range("K82").Value = 11838
$ws.Range("L82").Value = 24997.143
$ws.Range("M82").Value = -11455
$ws.Range("N82").Value = -25763.143
$ws.Range("H85").Value = 21049.4
$ws.Range("I85").Value = 11838
$ws.Range("J85").Value = 24997.143
$ws.Range("K85").Value = 11838
$ws.Range("L85").Value = 24997.143
$ws.Range("M85").Value = -10512
$ws.Range("N85").Value = -27649.143
$ws.Range("H88").Value = 56000
$ws.Range("J88").Value = 56000
$ws.Range("L88").Value = 56000
$ws.Range("N88").Value = -56812
$ws.Range("H91").Value = 56000
$ws.Range("J91").Value = 56000
$ws.Range("L91").Value = 56000
$ws.Range("N91").Value = -58808
$ws.Range("H94").Value = 1102.3334
$ws.Range("I94").Value = 1102.3334
$ws.Range("K94").Value = 1102.3334
$ws.Range("M94").Value = -651.3334
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 343.86487
$ws.Range("I7").Value = 181.80952
$ws.Range("J7").Value = 556.5625
$ws.Range("K7").Value = 181.80952
$ws.Range("L7").Value = 556.5625
$ws.Range("M7").Value = -68.80951999999999
$ws.Range("N7").Value = -782.5625
$ws.Range("H57").Value = 90000
$ws.Range("J57").Value = 100000
$ws.Range("L57").Value = 100000
$ws.Range("N57").Value = -101120
$ws.Range("H58").Value = 1566.1111
$ws.Range("I58").Value = 1349.6666
$ws.Range("K58").Value = 1349.6666
$ws.Range("M58").Value = -1146.6666
$ws.Range("H59").Value = 15907.091
$ws.Range("J59").Value = 15907.091
$ws.Range("L59").Value = 15907.091
$ws.Range("N59").Value = -18197.091
$ws.Range("H107").Value = 1681.2778
$ws.Range("I107").Value = 1130.8572
$ws.Range("J107").Value = 2031.5454
$ws.Range("K107").Value = 1130.8572
$ws.Range("L107").Value = 2031.5454
$ws.Range("M107").Value = 789.1428000000001
$ws.Range("N107").Value = -5871.5454
$ws.Range("H132").Value = 2904.5
$ws.Range("I132").Value = 2573.6924
$ws.Range("K132").Value = 7721.0772
$ws.Range("M132").Value = -5191.0772
$ws.Range("H134").Value = 2998
$ws.Range("I134").Value = 2998
$ws.Range("K134").Value = 8994
$ws.Range("M134").Value = -6459
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 1566.1111
$ws.Range("I136").Value = 1349.6666
$ws.Range("K136").Value = 4048.9998
$ws.Range("M136").Value = -1498.9998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 8999.799999999999
$ws.Range("J42").Value = 8999.799999999999
$ws.Range("L42").Value = 26999.4
$ws.Range("N42").Value = -28067.4
$ws.Range("H55").Value = 12602074
$ws.Range("I55").Value = 167206
$ws.Range("J55").Value = 31254376
$ws.Range("K55").Value = 501618
$ws.Range("L55").Value = 93763128
$ws.Range("M55").Value = -501441
$ws.Range("N55").Value = -93763482
$ws.Range("H68").Value = 1978.5
$ws.Range("I68").Value = 1665.6666
$ws.Range("J68").Value = 2166.2
$ws.Range("K68").Value = 4996.9998
$ws.Range("L68").Value = 6498.599999999999
$ws.Range("M68").Value = -4185.9998
$ws.Range("N68").Value = -8120.599999999999
$ws.Range("H71").Value = 1978.5
$ws.Range("I71").Value = 1665.6666
$ws.Range("J71").Value = 2166.2
$ws.Range("K71").Value = 14990.9994
$ws.Range("L71").Value = 19495.8
$ws.Range("M71").Value = -10934.9994
$ws.Range("N71").Value = -27607.8
$ws.Range("H107").Value = 613.9524
$ws.Range("J107").Value = 700.25
$ws.Range("L107").Value = 2100.75
$ws.Range("N107").Value = -5940.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 17746.871
$ws.Range("I97").Value = 24003.715
$ws.Range("J97").Value = 1820.3636
$ws.Range("K97").Value = 24003.715
$ws.Range("L97").Value = 1820.3636
$ws.Range("M97").Value = -23507.715
$ws.Range("N97").Value = -2812.3636
$ws.Range("H113").Value = 26318038
$ws.Range("I113").Value = 35716264
$ws.Range("J113").Value = 2998.8
$ws.Range("K113").Value = 35716264
$ws.Range("L113").Value = 2998.8
$ws.Range("M113").Value = -35714094
$ws.Range("N113").Value = -7338.8
$ws.Range("H126").Value = 3258.25
$ws.Range("I126").Value = 3182.5
$ws.Range("K126").Value = 9547.5
$ws.Range("M126").Value = -7077.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 499.875
$ws.Range("I16").Value = 413
$ws.Range("J16").Value = 586.75
$ws.Range("K16").Value = 413
$ws.Range("L16").Value = 586.75
$ws.Range("M16").Value = -243
$ws.Range("N16").Value = -926.75
$ws.Range("I55").Value = 364.46155
$ws.Range("J55").Value = 1649.6666
$ws.Range("K55").Value = 364.46155
$ws.Range("L55").Value = 1649.6666
$ws.Range("M55").Value = -191.46155
$ws.Range("N55").Value = -1995.6666
$ws.Range("H93").Value = 1960.8125
$ws.Range("I93").Value = 1683.8572
$ws.Range("K93").Value = 1683.8572
$ws.Range("M93").Value = -435.8571999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 23554.8
$ws.Range("J74").Value = 23554.8
$ws.Range("L74").Value = 23554.8
$ws.Range("N74").Value = -25426.8
$ws.Range("H77").Value = 23554.8
$ws.Range("J77").Value = 23554.8
$ws.Range("L77").Value = 70664.39999999999
$ws.Range("N77").Value = -80024.39999999999
$ws.Range("H136").Value = 872.5238000000001
$ws.Range("J136").Value = 799
$ws.Range("L136").Value = 2397
$ws.Range("N136").Value = -7497
